$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M53").Value = -971.7141999999999
$ws.Range("H53").Value = 812.3333
$ws.Range("N53").Value = -1389.5
$ws.Range("L53").Value = 115.5
$ws.Range("J53").Value = 115.5
$ws.Range("I53").Value = 1608.7142
$ws.Range("K53").Value = 1608.7142
$ws.Range("M54").Value = -69
$ws.Range("H54").Value = 555
$ws.Range("I54").Value = 555
$ws.Range("K54").Value = 555
$ws.Range("L55").Value = 196.66667
$ws.Range("H55").Value = 206
$ws.Range("K55").Value = 220
$ws.Range("N55").Value = -624.6666700000001
$ws.Range("J55").Value = 196.66667
$ws.Range("M55").Value = -6
$ws.Range("I55").Value = 220
$ws.Range("J100").Value = 2126.5
$ws.Range("M100").Value = -861.5
$ws.Range("K100").Value = 1402.5
$ws.Range("N100").Value = -3208.5
$ws.Range("H100").Value = 1981.7
$ws.Range("I100").Value = 1402.5
$ws.Range("L100").Value = 2126.5
$ws.Range("H129").Value = 2226.2559
$ws.Range("J129").Value = 2859.375
$ws.Range("L129").Value = 8578.125
$ws.Range("N129").Value = -18578.125
$ws.Range("L137").Value = 50003916
$ws.Range("H137").Value = 10527437
$ws.Range("I137").Value = 804.4286
$ws.Range("M137").Value = 136.7142000000003
$ws.Range("K137").Value = 2413.2858
$ws.Range("N137").Value = -50009016
$ws.Range("J137").Value = 16667972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L32").Value = 8959.117
$ws.Range("N32").Value = -9533.117
$ws.Range("I32").Value = 10622.667
$ws.Range("M32").Value = -10335.667
$ws.Range("K32").Value = 10622.667
$ws.Range("H32").Value = 10255.39
$ws.Range("J32").Value = 8959.117
$ws.Range("I45").Value = 2039.6364
$ws.Range("K45").Value = 2039.6364
$ws.Range("M45").Value = -1662.6364
$ws.Range("H45").Value = 2039.2858
$ws.Range("M53").Value = -3863
$ws.Range("H53").Value = 4545
$ws.Range("N53").ClearContents()
$ws.Range("L53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("I53").Value = 4545
$ws.Range("K53").Value = 4545
$ws.Range("J54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M61").Value = -8065405
$ws.Range("L61").Value = 1161.5
$ws.Range("I61").Value = 8065617
$ws.Range("H61").Value = 6945554
$ws.Range("K61").Value = 8065617
$ws.Range("J61").Value = 1161.5
$ws.Range("N61").Value = -1585.5
$ws.Range("J136").Value = 1161.5
$ws.Range("H136").Value = 6945554
$ws.Range("M136").Value = -24194301
$ws.Range("K136").Value = 24196851
$ws.Range("I136").Value = 8065617
$ws.Range("L136").Value = 3484.5
$ws.Range("N136").Value = -8584.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M54").Value = -1271.3334
$ws.Range("H54").Value = 1755.3334
$ws.Range("I54").Value = 1755.3334
$ws.Range("K54").Value = 1755.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M33").Value = -1366.5
$ws.Range("N33").Value = -5758
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 1745.5
$ws.Range("H33").Value = 2830.3333
$ws.Range("I33").Value = 1745.5
$ws.Range("L33").Value = 5000
$ws.Range("L58").Value = 2478.8572
$ws.Range("J58").Value = 2478.8572
$ws.Range("M58").Value = -986.5714
$ws.Range("N58").Value = -2884.8572
$ws.Range("I58").Value = 1189.5714
$ws.Range("K58").Value = 1189.5714
$ws.Range("H58").Value = 1705.2858
$ws.Range("J105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("K105").Value = 1199.8
$ws.Range("M105").Value = 547.2
$ws.Range("H105").Value = 1199.8
$ws.Range("I105").Value = 1199.8
$ws.Range("L105").Value = 0
$ws.Range("J136").Value = 2478.8572
$ws.Range("H136").Value = 1705.2858
$ws.Range("M136").Value = -1018.7142
$ws.Range("K136").Value = 3568.7142
$ws.Range("I136").Value = 1189.5714
$ws.Range("L136").Value = 7436.571599999999
$ws.Range("N136").Value = -12536.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 694.3333
$ws.Range("M47").Value = -1493.5
$ws.Range("J47").Value = 800
$ws.Range("I47").Value = 641.5
$ws.Range("N47").Value = -3262
$ws.Range("K47").Value = 1924.5
$ws.Range("L47").Value = 2400
$ws.Range("I49").Value = 2000
$ws.Range("M49").Value = -5844
$ws.Range("H49").Value = 4400
$ws.Range("K49").Value = 6000
$ws.Range("J54").Value = 2933.3333
$ws.Range("N54").Value = -9917.999899999999
$ws.Range("H54").Value = 2933.3333
$ws.Range("L54").Value = 8799.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H53").Value = 12000
$ws.Range("N53").Value = -13262
$ws.Range("L53").Value = 12000
$ws.Range("J53").Value = 12000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("H55").Value = 5150
$ws.Range("K55").Value = 5150
$ws.Range("N55").ClearContents()
$ws.Range("J55").Value = 0
$ws.Range("M55").Value = -4823
$ws.Range("I55").Value = 5150
$ws.Range("L138").Value = 58950
$ws.Range("H138").Value = 58950
$ws.Range("J138").Value = 58950
$ws.Range("N138").Value = -69230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L22").Value = 1171.6666
$ws.Range("K22").Value = 396.25
$ws.Range("H22").Value = 957.7586
$ws.Range("M22").Value = -101.25
$ws.Range("J22").Value = 1171.6666
$ws.Range("N22").Value = -1761.6666
$ws.Range("I22").Value = 396.25
$ws.Range("I27").Value = 396.25
$ws.Range("N27").Value = -1385.6666
$ws.Range("H27").Value = 957.7586
$ws.Range("L27").Value = 1171.6666
$ws.Range("J27").Value = 1171.6666
$ws.Range("M27").Value = -289.25
$ws.Range("K27").Value = 396.25
$ws.Range("K42").Value = 8000
$ws.Range("J42").Value = 0
$ws.Range("M42").Value = -7437
$ws.Range("N42").ClearContents()
$ws.Range("L42").Value = 0
$ws.Range("I42").Value = 8000
$ws.Range("H42").Value = 8000
$ws.Range("M46").Value = -415.3333
$ws.Range("H46").Value = 845.3333
$ws.Range("K46").Value = 603.3333
$ws.Range("I46").Value = 603.3333
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("L47").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("I49").Value = 8000
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("M49").Value = -7853
$ws.Range("H49").Value = 8000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 8000
$ws.Range("J52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("L52").Value = 0
$ws.Range("J54").Value = 33753.6
$ws.Range("N54").Value = -35041.6
$ws.Range("H54").Value = 33753.6
$ws.Range("L54").Value = 33753.6
$ws.Range("L55").Value = 996.6667
$ws.Range("H55").Value = 411.75
$ws.Range("K55").Value = 276.76923
$ws.Range("N55").Value = -1342.6667
$ws.Range("J55").Value = 996.6667
$ws.Range("M55").Value = -103.76923
$ws.Range("I55").Value = 276.76923
$ws.Range("N99").ClearContents()
$ws.Range("H99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("J99").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N48").ClearContents()
$ws.Range("J48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("L49").Value = 17708
$ws.Range("N49").Value = -18168
$ws.Range("M49").ClearContents()
$ws.Range("H49").Value = 17708
$ws.Range("J49").Value = 17708
$ws.Range("K49").Value = 0
$ws.Range("K132").Value = 7520.700000000001
$ws.Range("M132").Value = -4990.700000000001
$ws.Range("H132").Value = 2462.8076
$ws.Range("I132").Value = 2506.9
